# Auto-generated edit script applying numeric corrections to tradeskill profit sheets.
# Each sheet (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) gets a batch of cell updates
# to columns H-N reflecting refreshed market-board pricing data.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1430.125
$ws.Range("I19").Value = 1379.909
$ws.Range("J19").Value = 1540.6
$ws.Range("K19").Value = 1379.909
$ws.Range("L19").Value = 1540.6
$ws.Range("M19").Value = -1204.909
$ws.Range("N19").Value = -1890.6
$ws.Range("H34").Value = 5029.3335
$ws.Range("I34").Value = 5029.3335
$ws.Range("K34").Value = 5029.3335
$ws.Range("M34").Value = -4826.3335
$ws.Range("H36").Value = 5029.3335
$ws.Range("I36").Value = 5029.3335
$ws.Range("K36").Value = 5029.3335
$ws.Range("M36").Value = -4314.3335
$ws.Range("H70").Value = 1574.25
$ws.Range("J70").Value = 1497.5
$ws.Range("L70").Value = 4492.5
$ws.Range("N70").Value = -5032.5
$ws.Range("H73").Value = 1574.25
$ws.Range("J73").Value = 1497.5
$ws.Range("L73").Value = 4492.5
$ws.Range("N73").Value = -6364.5
$ws.Range("H106").Value = 6189
$ws.Range("I106").Value = 5985
$ws.Range("J106").Value = 6495
$ws.Range("K106").Value = 5985
$ws.Range("L106").Value = 6495
$ws.Range("M106").Value = -5354
$ws.Range("N106").Value = -7757
$ws.Range("H112").Value = 1372.4286
$ws.Range("J112").Value = 1372.4286
$ws.Range("L112").Value = 4117.2858
$ws.Range("N112").Value = -6333.2858
$ws.Range("H131").Value = 9750.25
$ws.Range("I131").Value = 9965.333000000001
$ws.Range("J131").Value = 9105
$ws.Range("K131").Value = 29895.999
$ws.Range("L131").Value = 27315
$ws.Range("M131").Value = -24855.999
$ws.Range("N131").Value = -37395
$ws.Range("H137").Value = 3615.75
$ws.Range("I137").Value = 3178.6
$ws.Range("J137").Value = 3928
$ws.Range("K137").Value = 9535.799999999999
$ws.Range("L137").Value = 11784
$ws.Range("M137").Value = -6985.799999999999
$ws.Range("N137").Value = -16884
$ws.Range("H138").Value = 8185.25
$ws.Range("J138").Value = 8185.25
$ws.Range("L138").Value = 24555.75
$ws.Range("N138").Value = -34835.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3766.577
$ws.Range("I32").Value = 3289.25
$ws.Range("K32").Value = 3289.25
$ws.Range("M32").Value = -3002.25
$ws.Range("H74").Value = 2063.3684
$ws.Range("I74").Value = 685.1539
$ws.Range("J74").Value = 5049.5
$ws.Range("K74").Value = 685.1539
$ws.Range("L74").Value = 5049.5
$ws.Range("M74").Value = 188.8461
$ws.Range("N74").Value = -6797.5
$ws.Range("H77").Value = 2063.3684
$ws.Range("I77").Value = 685.1539
$ws.Range("J77").Value = 5049.5
$ws.Range("K77").Value = 3425.7695
$ws.Range("L77").Value = 25247.5
$ws.Range("M77").Value = 942.2304999999997
$ws.Range("N77").Value = -33983.5
$ws.Range("H97").Value = 752.8889
$ws.Range("I97").Value = 860.8570999999999
$ws.Range("K97").Value = 860.8570999999999
$ws.Range("M97").Value = -364.8570999999999
$ws.Range("H102").Value = 1410.875
$ws.Range("I102").Value = 1556.3334
$ws.Range("K102").Value = 1556.3334
$ws.Range("M102").Value = 65.66660000000002
$ws.Range("H122").Value = 1531.4615
$ws.Range("I122").Value = 1531.4615
$ws.Range("K122").Value = 4594.3845
$ws.Range("M122").Value = -2144.3845

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 31761
$ws.Range("I82").Value = 12500
$ws.Range("K82").Value = 12500
$ws.Range("M82").Value = -12117
$ws.Range("H85").Value = 31761
$ws.Range("I85").Value = 12500
$ws.Range("K85").Value = 12500
$ws.Range("M85").Value = -11174

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 4999
$ws.Range("J10").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("N10").ClearContents()
$ws.Range("H16").Value = 3802
$ws.Range("J16").Value = 4329.6665
$ws.Range("L16").Value = 4329.6665
$ws.Range("N16").Value = -4903.6665
$ws.Range("H31").Value = 4144.879
$ws.Range("I31").Value = 1310.3334
$ws.Range("J31").Value = 7546.3335
$ws.Range("K31").Value = 1310.3334
$ws.Range("L31").Value = 7546.3335
$ws.Range("M31").Value = -1015.3334
$ws.Range("N31").Value = -8136.3335
$ws.Range("H34").Value = 4144.879
$ws.Range("I34").Value = 1310.3334
$ws.Range("J34").Value = 7546.3335
$ws.Range("K34").Value = 1310.3334
$ws.Range("L34").Value = 7546.3335
$ws.Range("M34").Value = -1108.3334
$ws.Range("N34").Value = -7950.3335
$ws.Range("H113").Value = 3802
$ws.Range("J113").Value = 4329.6665
$ws.Range("L113").Value = 4329.6665
$ws.Range("N113").Value = -8669.666499999999
$ws.Range("H122").Value = 1724.25
$ws.Range("I122").Value = 1724.25
$ws.Range("K122").Value = 5172.75
$ws.Range("M122").Value = -2722.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 2000
$ws.Range("J39").Value = 2000
$ws.Range("L39").Value = 6000
$ws.Range("N39").Value = -6588
$ws.Range("H137").Value = 21199.8
$ws.Range("J137").Value = 19249.75
$ws.Range("L137").Value = 57749.25
$ws.Range("N137").Value = -67949.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1061.3
$ws.Range("I97").Value = 777.1429000000001
$ws.Range("J97").Value = 1724.3334
$ws.Range("K97").Value = 777.1429000000001
$ws.Range("L97").Value = 1724.3334
$ws.Range("M97").Value = -281.1429000000001
$ws.Range("N97").Value = -2716.3334
$ws.Range("H102").Value = 3945.75
$ws.Range("I102").Value = 3080.8572
$ws.Range("K102").Value = 3080.8572
$ws.Range("M102").Value = -1458.8572
$ws.Range("H132").Value = 3321.25
$ws.Range("I132").Value = 2844.75
$ws.Range("J132").Value = 4750.75
$ws.Range("K132").Value = 8534.25
$ws.Range("L132").Value = 14252.25
$ws.Range("M132").Value = -6004.25
$ws.Range("N132").Value = -19312.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 890
$ws.Range("I7").Value = 890
$ws.Range("K7").Value = 890
$ws.Range("M7").Value = -778
$ws.Range("H40").Value = 3166.0527
$ws.Range("I40").Value = 3126.8572
$ws.Range("J40").Value = 3275.8
$ws.Range("K40").Value = 3126.8572
$ws.Range("L40").Value = 3275.8
$ws.Range("M40").Value = -2990.8572
$ws.Range("N40").Value = -3547.8
$ws.Range("H68").Value = 2741.3333
$ws.Range("I68").Value = 2808.4546
$ws.Range("J68").Value = 2003
$ws.Range("K68").Value = 2808.4546
$ws.Range("L68").Value = 2003
$ws.Range("M68").Value = -2059.4546
$ws.Range("N68").Value = -3501
$ws.Range("H71").Value = 2741.3333
$ws.Range("I71").Value = 2808.4546
$ws.Range("J71").Value = 2003
$ws.Range("K71").Value = 14042.273
$ws.Range("L71").Value = 10015
$ws.Range("M71").Value = -10298.273
$ws.Range("N71").Value = -17503
$ws.Range("H126").Value = 890
$ws.Range("I126").Value = 890
$ws.Range("K126").Value = 2670
$ws.Range("M126").Value = -200
$ws.Range("H136").Value = 34764
$ws.Range("I136").Value = 7899.5713
$ws.Range("K136").Value = 23698.7139
$ws.Range("M136").Value = -21148.7139
$ws.Range("H141").Value = 114899
$ws.Range("J141").Value = 114899
$ws.Range("L141").Value = 114899
$ws.Range("N141").Value = -125259

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 3320.7778
$ws.Range("I96").Value = 3485.875
$ws.Range("J96").Value = 2000
$ws.Range("K96").Value = 3485.875
$ws.Range("L96").Value = 2000
$ws.Range("M96").Value = -2112.875
$ws.Range("N96").Value = -4746
$ws.Range("H122").Value = 1581.25
$ws.Range("I122").Value = 1420.0667
$ws.Range("K122").Value = 4260.2001
$ws.Range("M122").Value = -1810.2001

